# ---------------------------------------------------------------------------
# Refresh the NATMI ligand-receptor (Ntf5 -> Ntrk1) edge table with values
# recomputed from the new TPM expression matrix.
#
# The sending/target cluster universe changes from
#   {FAPs, Inflammatory-Mac, Neutrophils, Resolving-Mac} x {FAPs, MuSCs}   (8 rows)
# to
#   {ECs, FAPs, Neutrophils, Resolving-Mac} x {ECs, FAPs, MuSCs}           (12 rows)
#
# i.e. "Inflammatory-Mac" is replaced by "ECs" as a sending cluster, and "ECs"
# is added as a third possible target cluster alongside "FAPs"/"MuSCs". All
# numeric expression/specificity columns (E:T) are refreshed to match.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs (sending) -> ECs (target)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf5"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.2515545
$ws.Range("H2").Value = 0.503109
$ws.Range("I2").Value = 0.3771142493566836
$ws.Range("J2").Value = 0.2875567487178491
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0003185
$ws.Range("N2").Value = 0.000637
$ws.Range("O2").Value = 0.001593934541086978
$ws.Range("P2").Value = 0.001183922444716212
$ws.Range("Q2").Value = [double]"8.012010825E-05"
$ws.Range("R2").Value = 0.000320480433
$ws.Range("S2").Value = 0.0006010954279857057
$ws.Range("T2").Value = 0.0003404448889366813

# Row 3: ECs (sending) -> FAPs (target)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf5"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.2515545
$ws.Range("H3").Value = 0.503109
$ws.Range("I3").Value = 0.3771142493566836
$ws.Range("J3").Value = 0.2875567487178491
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.138402
$ws.Range("N3").Value = 0.415206
$ws.Range("O3").Value = 0.6926333700330297
$ws.Range("P3").Value = 0.7716981202210981
$ws.Range("Q3").Value = 0.034815645909
$ws.Range("R3").Value = 0.208893875454
$ws.Range("S3").Value = 0.2612019134193961
$ws.Range("T3").Value = 0.2219070024424548

# Row 4: ECs (sending) -> MuSCs (target)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf5"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.2515545
$ws.Range("H4").Value = 0.503109
$ws.Range("I4").Value = 0.3771142493566836
$ws.Range("J4").Value = 0.2875567487178491
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.0610995
$ws.Range("N4").Value = 0.122199
$ws.Range("O4").Value = 0.3057726954258833
$ws.Range("P4").Value = 0.2271179573341859
$ws.Range("Q4").Value = 0.01536985417275
$ws.Range("R4").Value = 0.061479416691
$ws.Range("S4").Value = 0.1153112405093018
$ws.Range("T4").Value = 0.06530930138645766

# Row 5: FAPs (sending) -> ECs (target)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf5"
$ws.Range("C5").Value = "Ntrk1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.121806
$ws.Range("H5").Value = 0.365418
$ws.Range("I5").Value = 0.182603683325642
$ws.Range("J5").Value = 0.2088581440661546
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0003185
$ws.Range("N5").Value = 0.000637
$ws.Range("O5").Value = 0.001593934541086978
$ws.Range("P5").Value = 0.001183922444716212
$ws.Range("Q5").Value = [double]"3.8795211E-05"
$ws.Range("R5").Value = 0.000232771266
$ws.Range("S5").Value = 0.000291058318182449
$ws.Range("T5").Value = 0.0002472718445216926

# Row 6: FAPs (sending) -> FAPs (target)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf5"
$ws.Range("C6").Value = "Ntrk1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.121806
$ws.Range("H6").Value = 0.365418
$ws.Range("I6").Value = 0.182603683325642
$ws.Range("J6").Value = 0.2088581440661546
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.138402
$ws.Range("N6").Value = 0.415206
$ws.Range("O6").Value = 0.6926333700330297
$ws.Range("P6").Value = 0.7716981202210981
$ws.Range("Q6").Value = 0.016858194012
$ws.Range("R6").Value = 0.151723746108
$ws.Range("S6").Value = 0.1264774045622836
$ws.Range("T6").Value = 0.1611754371687188

# Row 7: FAPs (sending) -> MuSCs (target)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf5"
$ws.Range("C7").Value = "Ntrk1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.121806
$ws.Range("H7").Value = 0.365418
$ws.Range("I7").Value = 0.182603683325642
$ws.Range("J7").Value = 0.2088581440661546
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.0610995
$ws.Range("N7").Value = 0.122199
$ws.Range("O7").Value = 0.3057726954258833
$ws.Range("P7").Value = 0.2271179573341859
$ws.Range("Q7").Value = 0.007442285697000001
$ws.Range("R7").Value = 0.044653714182
$ws.Range("S7").Value = 0.05583522044517597
$ws.Range("T7").Value = 0.04743543505291414

# Row 8: Neutrophils (sending) -> ECs (target)
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Ntf5"
$ws.Range("C8").Value = "Ntrk1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.244596
$ws.Range("H8").Value = 0.733788
$ws.Range("I8").Value = 0.3666825158589784
$ws.Range("J8").Value = 0.4194035318950228
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.0003185
$ws.Range("N8").Value = 0.000637
$ws.Range("O8").Value = 0.001593934541086978
$ws.Range("P8").Value = 0.001183922444716212
$ws.Range("Q8").Value = [double]"7.7903826E-05"
$ws.Range("R8").Value = 0.000467422956
$ws.Range("S8").Value = 0.0005844679276402993
$ws.Range("T8").Value = 0.0004965412548037692

# Row 9: Neutrophils (sending) -> FAPs (target)
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Ntf5"
$ws.Range("C9").Value = "Ntrk1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.244596
$ws.Range("H9").Value = 0.733788
$ws.Range("I9").Value = 0.3666825158589784
$ws.Range("J9").Value = 0.4194035318950228
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.138402
$ws.Range("N9").Value = 0.415206
$ws.Range("O9").Value = 0.6926333700330297
$ws.Range("P9").Value = 0.7716981202210981
$ws.Range("Q9").Value = 0.033852575592
$ws.Range("R9").Value = 0.3046731803279999
$ws.Range("S9").Value = 0.253976546691594
$ws.Range("T9").Value = 0.3236529171774785

# Row 10: Neutrophils (sending) -> MuSCs (target)
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Ntf5"
$ws.Range("C10").Value = "Ntrk1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.244596
$ws.Range("H10").Value = 0.733788
$ws.Range("I10").Value = 0.3666825158589784
$ws.Range("J10").Value = 0.4194035318950228
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.0610995
$ws.Range("N10").Value = 0.122199
$ws.Range("O10").Value = 0.3057726954258833
$ws.Range("P10").Value = 0.2271179573341859
$ws.Range("Q10").Value = 0.014944693302
$ws.Range("R10").Value = 0.08966815981200001
$ws.Range("S10").Value = 0.112121501239744
$ws.Range("T10").Value = 0.09525407346274065

# Row 11: Resolving-Mac (sending) -> ECs (target)
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Ntf5"
$ws.Range("C11").Value = "Ntrk1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04909466666666667
$ws.Range("H11").Value = 0.147284
$ws.Range("I11").Value = 0.07359955145869621
$ws.Range("J11").Value = 0.08418157532097355
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.5
$ws.Range("M11").Value = 0.0003185
$ws.Range("N11").Value = 0.000637
$ws.Range("O11").Value = 0.001593934541086978
$ws.Range("P11").Value = 0.001183922444716212
$ws.Range("Q11").Value = [double]"1.563665133333333E-05"
$ws.Range("R11").Value = [double]"9.381990799999999E-05"
$ws.Range("S11").Value = 0.0001173128672785244
$ws.Range("T11").Value = [double]"9.966445645406895E-05"

# Row 12: Resolving-Mac (sending) -> FAPs (target)
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Ntf5"
$ws.Range("C12").Value = "Ntrk1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04909466666666667
$ws.Range("H12").Value = 0.147284
$ws.Range("I12").Value = 0.07359955145869621
$ws.Range("J12").Value = 0.08418157532097355
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.138402
$ws.Range("N12").Value = 0.415206
$ws.Range("O12").Value = 0.6926333700330297
$ws.Range("P12").Value = 0.7716981202210981
$ws.Range("Q12").Value = 0.006794800056
$ws.Range("R12").Value = 0.06115320050399999
$ws.Range("S12").Value = 0.05097750535975614
$ws.Range("T12").Value = 0.06496276343244607

# Row 13: Resolving-Mac (sending) -> MuSCs (target)
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Ntf5"
$ws.Range("C13").Value = "Ntrk1"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04909466666666667
$ws.Range("H13").Value = 0.147284
$ws.Range("I13").Value = 0.07359955145869621
$ws.Range("J13").Value = 0.08418157532097355
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.0610995
$ws.Range("N13").Value = 0.122199
$ws.Range("O13").Value = 0.3057726954258833
$ws.Range("P13").Value = 0.2271179573341859
$ws.Range("Q13").Value = 0.002999659586
$ws.Range("R13").Value = 0.017997957516
$ws.Range("S13").Value = 0.02250473323166154
$ws.Range("T13").Value = 0.01911914743207342
